# Apply the edits described by the commit diff to Sheet1 of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: re-apply one of the two existing conditional-color styles -----
# Style index 3 = light-green fill ("doğru / correct"),
# Style index 4 = light-red   fill ("yanlış / incorrect").
# Both styles already exist in the workbook; we copy the *format only* from a
# donor cell that is known to keep that style so that we reuse the existing
# style index instead of creating a brand new one.
$styleDonor3 = $ws.Range("E2")   # untouched cell carrying style index 3 (green)
$styleDonor4 = $ws.Range("G5")   # untouched cell carrying style index 4 (red)

function Set-CellStyle($cellRef, [int]$styleNumber) {
    if ($styleNumber -eq 3) {
        $styleDonor3.Copy() | Out-Null
    } else {
        $styleDonor4.Copy() | Out-Null
    }
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# --- Row 2 -------------------------------------------------------------
$ws.Range("I2").Value = "Siparişte eksik ürün bildirimi"

# --- Row 3 -------------------------------------------------------------
Set-CellStyle "G3" 3
$ws.Range("G3").Value = "Sorun"

# --- Row 6 -------------------------------------------------------------
Set-CellStyle "H6" 3
$ws.Range("H6").Value = "İptal"
$ws.Range("I6").Value = "Sipariş iptali"

# --- Row 8 -------------------------------------------------------------
$ws.Range("H8").Value = "Hasarlı ürün"
$ws.Range("I8").Value = "Defolu ürün için iade/değişim talebi, fotoğraf kanıtı gerektiği konusunda memnuniyetsizlik"

# --- Row 9 -------------------------------------------------------------
$ws.Range("I9").Value = "Üyelik bilgilerimi güncellemek"

# --- Row 10 ------------------------------------------------------------
$ws.Range("I10").Value = "Fatura üzerinde hatalı fiyatlar tespiti ve düzeltme talebi"

# --- Row 11 ------------------------------------------------------------
$ws.Range("I11").Value = "İptal butonunun görünmemesi ve iptal edilememesi nedeniyle hizmetten memnuniyetsizlik"

# --- Row 12 ------------------------------------------------------------
$ws.Range("I12").Value = "Ürün stok durumu ve stok gelecekteki geliş tarihi bilgisi"

# --- Row 13 ------------------------------------------------------------
Set-CellStyle "G13" 4
$ws.Range("G13").Value = "Sorun"
Set-CellStyle "H13" 3
$ws.Range("H13").Value = "Ödeme"
$ws.Range("I13").Value = "Kredi kartı iki kez çekim işlemi (çekim hatası)"

# --- Row 15 ------------------------------------------------------------
Set-CellStyle "E15" 4
$ws.Range("E15").Value = "Çözüldü"
$ws.Range("I15").Value = "Kredi kartı reddediliyor"

# --- Row 16 ------------------------------------------------------------
$ws.Range("I16").Value = "Ürün beden tablosu bilgisini bulma"

# --- Row 17 ------------------------------------------------------------
$ws.Range("I17").Value = "Adres yanlış girdim ve düzeltme yolu soruldu fakat çözüm sağlayamadı"

# --- Row 19 ------------------------------------------------------------
Set-CellStyle "E19" 3
$ws.Range("E19").Value = "Çözüldü"
Set-CellStyle "F19" 3
$ws.Range("F19").Value = "Nötr"
$ws.Range("I19").Value = "Web sitesi hatası (404 Not Found) ile ilgili bilgi talebi/şikayet"

# --- Row 20 ------------------------------------------------------------
$ws.Range("I20").Value = "Hesap üzerinden çekim yapılmış ancak sipariş görünmüyor/ödemeye dair sorun"

# --- Row 21 ------------------------------------------------------------
Set-CellStyle "F21" 3
$ws.Range("F21").Value = "Nötr"
$ws.Range("I21").Value = "Ürün özelliği hakkında bilgi (kumaş niteliği)"

# --- Row 22 ------------------------------------------------------------
$ws.Range("I22").Value = "Siparişte eksik/promosyon ürünlerin olması ve çözüm talebi"

# --- Row 23 ------------------------------------------------------------
$ws.Range("I23").Value = "İade süreci gecikmesi ve para iadesi ile ilgili memnuniyetsizlik"

# --- Row 24 ------------------------------------------------------------
$ws.Range("I24").Value = "Kargo kaynaklı hasar bildirimi ve değişim talebi"

# --- Row 25 ------------------------------------------------------------
$ws.Range("I25").Value = "Değişim seçeneğinin olmaması ve iade+yeniden siparişin zahmetli oluşu"

# --- Row 27 ------------------------------------------------------------
$ws.Range("I27").Value = "Yanlış ürün gönderimi nedeniyle iade/yeniden gönderim süreci başlatıldı"

# --- Row 28 ------------------------------------------------------------
Set-CellStyle "F28" 3
$ws.Range("F28").Value = "Nötr"

# --- Row 29 ------------------------------------------------------------
$ws.Range("I29").Value = "Siparişte eksik ürün ve stok hatası nedeniyle sorun yaşanıyor"

# --- Row 30 ------------------------------------------------------------
$ws.Range("G30").Value = "İade"
$ws.Range("I30").Value = "İade süreciyle ilgili sorun çözümü ve yeni iade kodu oluşturuldu"

# --- Row 31 ------------------------------------------------------------
$ws.Range("I31").Value = "Teslimat adresinin yanlış yazılması ve düzeltilmesi"

# --- Row 33 ------------------------------------------------------------
$ws.Range("I33").Value = "Hesap üzerinden ödeme çekilmemesi/tekrar deneme?"

# --- Row 34 ------------------------------------------------------------
$ws.Range("I34").Value = "Beden uyuşmazlığı nedeniyle iade"

# --- Row 36 ------------------------------------------------------------
$ws.Range("G36").Value = "Sorun"
$ws.Range("I36").Value = "Kargo takibi ve gecikme"

# --- Row 37 ------------------------------------------------------------
Set-CellStyle "F37" 3
$ws.Range("F37").Value = "Pozitif"

# --- Row 38 ------------------------------------------------------------
$ws.Range("H38").Value = "Defolu ürün"
$ws.Range("I38").Value = "Kulaklıktaki arıza ve garanti kapsamında tamir talebi"

# --- Row 39 ------------------------------------------------------------
Set-CellStyle "F39" 3
$ws.Range("F39").Value = "Nötr"
$ws.Range("I39").Value = "Abonelikten çıkma talebi"

# --- Row 40 ------------------------------------------------------------
$ws.Range("I40").Value = "Kargonun gecikmesi nedeniyle teslimatın gecikmesi"

# --- Row 41 ------------------------------------------------------------
$ws.Range("I41").Value = "İade süreci uzaması ve kötü hizmet"

# --- Summary rows --------------------------------------------------------
$ws.Range("B46").Value = 28
$ws.Range("B47").Value = 12
$ws.Range("B48").Value = "%70.0"
